# lesson 499 - vocab and homework
#
# 1. Fix two typos in the vocabulary list (paragraph near the top of the
#    document): "sparsley" -> "sparsely", "decreptitude" -> "decrepitude".
# 2. Fill in seven of the "....." answer blanks throughout the reading
#    passage with the missing vocabulary word (splitting the existing run
#    of ellipsis characters after its first character, exactly as the
#    author typed the answer in the middle of the blank).
# 3. Move the "_GoBack" bookmark from the blank after "taking" (paragraph
#    about "A growing number of students ...") to sit right after the
#    newly typed word "default" in the blank following "by ...".

function Insert-WordAfterNeedle {
    param($doc, $needle, $wordToInsert)

    $t = $doc.Content.Text
    $idx = $t.IndexOf($needle)
    if ($idx -lt 0) {
        Write-Output "NOT FOUND: [$needle]"
        return -1
    }
    $insertPos = $idx + $needle.Length
    $rng = $doc.Range($insertPos, $insertPos)
    $rng.InsertAfter($wordToInsert)
    return $insertPos + $wordToInsert.Length
}

$d = $word.ActiveDocument

$ellipsis = [char]8230   # "…" - the single ellipsis glyph used throughout the blanks
$nbsp = [char]160        # non-breaking space that separates the blanks from surrounding text

# --- 1. Typo fixes in the vocabulary list -----------------------------------
[void]$d.Content.Find.Execute("sparsley", $true, $false, $false, $false, $false, $true, 1, $false, "sparsely", 2)
[void]$d.Content.Find.Execute("decreptitude", $true, $false, $false, $false, $false, $true, 1, $false, "decrepitude", 2)

# --- 2. Fill in the answer blanks -------------------------------------------

# "Stop going ……………………………" -> "Stop going …back and forth…………………………"
$needleGoing = "going " + $ellipsis
[void](Insert-WordAfterNeedle $d $needleGoing "back and forth")

# "He ……………………………. himself" -> "He …immersed…………………………. himself"
$needleHe = "He" + $nbsp + $ellipsis
[void](Insert-WordAfterNeedle $d $needleHe "immersed")

# "why I …………………………. support" -> "why I …wholeheartedly………………………. support"
$needleI = " I" + $nbsp + $ellipsis
[void](Insert-WordAfterNeedle $d $needleI "wholeheartedly")

# "have that ……………………. ability" -> "have that …latent………………….  ability"
$needleThat = "that" + $nbsp + $ellipsis
[void](Insert-WordAfterNeedle $d $needleThat "latent")

# "fact of his ………………………….." -> "fact of his …decreptitude……………………….."
$needleHis = "his" + $nbsp + $ellipsis
[void](Insert-WordAfterNeedle $d $needleHis "decreptitude")

# "smart by …………………..." -> "smart by …default………………..."  (note bookmark, handled below)
$needleBy = "by " + $ellipsis
$afterDefaultPos = Insert-WordAfterNeedle $d $needleBy "default"

# "end of his ………………………….." -> "end of his …sabbatical……………………….."
$needleEndOfHis = "end of " + "his" + $nbsp + $ellipsis
[void](Insert-WordAfterNeedle $d $needleEndOfHis "sabbatical")

# --- 3. Relocate the "_GoBack" bookmark -------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
if ($afterDefaultPos -ge 0) {
    $bmRange = $d.Range($afterDefaultPos, $afterDefaultPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "done"
